$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.979.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.38%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.052.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.13%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.18%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'245.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.659"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.18%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'58.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.47%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.20%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0772"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.65%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.111"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.30%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'15.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.70%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.892"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +8.67%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.350.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.24%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.996.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.98%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'18.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.62%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'36.939.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.44%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'73.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0885"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.39%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'238.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E24").Value = "'  +1.42%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.86%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'168.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.06%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.97%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'20.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +12.34%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.41%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.33%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0613"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.24%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.08%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.0858"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.31%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.87%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.43%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.33%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.49%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.39%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0223"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.32%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -10.83%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.75%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'97.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.19%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'16.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.54%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.299.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.68%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.43%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.07%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'6.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.42%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.236.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.25%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'44.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.43%  "
$ws.Range("E51").Style = "Normal"
